$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Update the "Blocked matrix and MG" coupled-MG table (rows 13-18)
#    with corrected numbers.
# -----------------------------------------------------------------
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = 1.07
$ws.Range("H13").Value = 8
$ws.Range("L13").Value = 8

$ws.Range("D14").Value = 18
$ws.Range("E14").Value = 2.12
$ws.Range("H14").Value = 18
$ws.Range("L14").Value = 16

$ws.Range("D15").Value = 32
$ws.Range("E15").Value = 3.58
$ws.Range("H15").Value = 32
$ws.Range("L15").Value = 30

$ws.Range("D16").Value = 71
$ws.Range("E16").Value = 7.82
$ws.Range("H16").Value = 72
$ws.Range("L16").Value = 76

$ws.Range("D17").Value = 270
$ws.Range("E17").Value = 31.2
$ws.Range("H17").Value = 339
$ws.Range("L17").Value = 341

# D18/H18 used to show the ">5000" label; that text moves to the new
# table heading below, so these now show ">1000" (same text as F6/B16 cells).
$ws.Range("D18").Value = ">1000"
$ws.Range("H18").Value = ">1000"
$ws.Range("L18").Value = 1517

# -----------------------------------------------------------------
# 2) Add the new "Coupled MG with Anton's favorite solver options"
#    table at rows 36-46 (a copy of the layout used at rows 8-18,
#    minus the extra viscosity-contrast columns I/M).
# -----------------------------------------------------------------

$ws.Range("A36").Value = "Coupled MG with Anton's favorite solver options"
$ws.Range("A8").Copy()
$ws.Range("A36").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A37").Value = "Setup: FallingBlock_canonical_coupledMG_directCoarse.dat and MultipleSpheres_canonical_coupledMG_directCoarse.dat"
$ws.Range("A9").Copy()
$ws.Range("A37").PasteSpecial(-4122)

$ws.Range("C38").Value = "COUPLED MG with -pcmat_no_dev_proj option"
$ws.Range("C10").Copy()
$ws.Range("C38").PasteSpecial(-4122)

$ws.Range("H38").Value = "COUPLED MG without -pcmat_no_dev_proj option"
$ws.Range("H10").Copy()
$ws.Range("H38").PasteSpecial(-4122)

$ws.Range("L38").Value = "COUPLED MG with -pcmat_no_dev_proj option but computing viscosity contrast by having a matrix viscosity 1 and block viscosity as indicated"
$ws.Range("L10").Copy()
$ws.Range("L38").PasteSpecial(-4122)

$ws.Range("B39").Value = "Falling Block Setup"
$ws.Range("B11").Copy()
$ws.Range("B39").PasteSpecial(-4122)

$ws.Range("C11").Copy()
$ws.Range("C39").PasteSpecial(-4122)

$ws.Range("D39").Value = "8 Spheres setup"
$ws.Range("D11").Copy()
$ws.Range("D39").PasteSpecial(-4122)

$ws.Range("E11").Copy()
$ws.Range("E39").PasteSpecial(-4122)

$ws.Range("H39").Value = "8 Spheres setup"
$ws.Range("H11").Copy()
$ws.Range("H39").PasteSpecial(-4122)

$ws.Range("L39").Value = "8 Spheres setup"
$ws.Range("L11").Copy()
$ws.Range("L39").PasteSpecial(-4122)

$ws.Range("A40").Value = "Viscosity contrast"
$ws.Range("A12").Copy()
$ws.Range("A40").PasteSpecial(-4122)

$ws.Range("B40").Value = "# outer KSP it"
$ws.Range("B12").Copy()
$ws.Range("B40").PasteSpecial(-4122)

$ws.Range("C40").Value = "Total solve [s]"
$ws.Range("C12").Copy()
$ws.Range("C40").PasteSpecial(-4122)

$ws.Range("D40").Value = "# outer KSP it"
$ws.Range("D12").Copy()
$ws.Range("D40").PasteSpecial(-4122)

$ws.Range("E40").Value = "Total solve [s]"
$ws.Range("E12").Copy()
$ws.Range("E40").PasteSpecial(-4122)

$ws.Range("H40").Value = "# outer KSP it"
$ws.Range("H12").Copy()
$ws.Range("H40").PasteSpecial(-4122)

$ws.Range("L40").Value = "# outer KSP it"
$ws.Range("L12").Copy()
$ws.Range("L40").PasteSpecial(-4122)

# Data rows 41-46 (viscosity contrast 1 .. 100000)
$ws.Range("A41").Value = 1
$ws.Range("A13").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("B41").Value = 7
$ws.Range("D41").Value = 7
$ws.Range("H41").Value = 7
$ws.Range("L41").Value = 7

$ws.Range("A42").Value = 10
$ws.Range("A14").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("B42").Value = 11
$ws.Range("D42").Value = 12
$ws.Range("H42").Value = 12
$ws.Range("L42").Value = 13

$ws.Range("A43").Value = 100
$ws.Range("A15").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$ws.Range("B43").Value = 12
$ws.Range("D43").Value = 19
$ws.Range("H43").Value = 19
$ws.Range("L43").Value = 23

$ws.Range("A44").Value = 1000
$ws.Range("A16").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("B44").Value = 13
$ws.Range("D44").Value = 39
$ws.Range("H44").Value = 41
$ws.Range("L44").Value = 48

$ws.Range("A45").Value = 10000
$ws.Range("A17").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("B45").Value = 17
$ws.Range("D45").Value = 80
$ws.Range("H45").Value = 141
$ws.Range("L45").Value = 79

$ws.Range("A46").Value = 100000
$ws.Range("A18").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("B46").Value = 26

$ws.Range("D46").Value = 88
$ws.Range("B28").Copy()
$ws.Range("D46").PasteSpecial(-4122)

$ws.Range("E46").Value = 5.51

$ws.Range("H46").Value = 585

$ws.Range("L46").Value = "25?"
$ws.Range("D18").Copy()
$ws.Range("L46").PasteSpecial(-4122)

# -----------------------------------------------------------------
# 3) Sheet-level view bookkeeping to match the new extent of data.
# -----------------------------------------------------------------
$ws.Range("A12").Select()
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("E47").Select()

Write-Output "done"
